$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, pushing existing rows 193:254 down to 194:255
$ws.Rows("193:193").Insert()

# Row 193 keeps the same dimension/lookup values as its neighbours (Vega Monumental
# Concepcion, Biobio, Ajo, Chino, Primera, $/caja 10 kilos, China, Hortaliza) but
# carries its own date/volume/price figures.
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 44985
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = 100112003
$ws.Range("G193").Value = "Ajo"
$ws.Range("H193").Value = "Chino"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 220
$ws.Range("K193").Value = 15000
$ws.Range("L193").Value = 16000
$ws.Range("M193").Value = 15455
$ws.Range("N193").Value = "$/caja 10 kilos"
$ws.Range("O193").Value = "China"
$ws.Range("P193").Value = 1546
$ws.Range("Q193").Value = 10
$ws.Range("R193").Value = "Hortaliza"
